$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-09-21 Thursday" "2023-09-22 Friday"

Replace-Text "32÷3=10, 2" "78÷3=26, 0"
Replace-Text "56÷6=9, 2" "88÷5=17, 3"
Replace-Text "60÷6=10, 0" "77÷9=8, 5"
Replace-Text "24÷7=3, 3" "23÷2=11, 1"
Replace-Text "64÷5=12, 4" "12÷3=4, 0"

Replace-Text "64÷6=10, 4" "68÷3=22, 2"
Replace-Text "27÷8=3, 3" "62÷9=6, 8"
Replace-Text "43÷8=5, 3" "61÷4=15, 1"
Replace-Text "20÷5=4, 0" "21÷7=3, 0"
Replace-Text "58÷5=11, 3" "19÷6=3, 1"

Replace-Text "96÷9=10, 6" "75÷8=9, 3"
Replace-Text "99÷6=16, 3" "60÷8=7, 4"
Replace-Text "56÷3=18, 2" "50÷9=5, 5"
Replace-Text "92÷2=46, 0" "69÷2=34, 1"
Replace-Text "45÷2=22, 1" "38÷7=5, 3"

Replace-Text "41÷5=8, 1" "67÷3=22, 1"
Replace-Text "36÷8=4, 4" "38÷7=5, 3"
Replace-Text "96÷7=13, 5" "10÷3=3, 1"
Replace-Text "46÷4=11, 2" "73÷3=24, 1"
Replace-Text "33÷6=5, 3" "16÷9=1, 7"

Replace-Text "89÷5=17, 4" "68÷6=11, 2"
Replace-Text "48÷3=16, 0" "66÷7=9, 3"
Replace-Text "87÷9=9, 6" "27÷6=4, 3"
Replace-Text "45÷7=6, 3" "34÷8=4, 2"
Replace-Text "95÷5=19, 0" "51÷8=6, 3"
